$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (AL1:AN1): new period labels, matching existing header style ---
$ws.Range("AL1").Value = "31/12/2023"
$ws.Range("AM1").Value = "31/03/2024"
$ws.Range("AN1").Value = "30/06/2024"
$ws.Range("AK1").Copy()
$ws.Range("AL1:AN1").PasteSpecial(-4122)

# --- Numeric data rows: append AL/AM/AN (cols 38-40) values for the three new periods ---
$data = @{
    2 = @(2543864.064, 2615521.024, 2532305.92)
    3 = @(573198.976, 651928, 533417.9840000001)
    4 = @(35566, 36996, 26724)
    5 = @(261052.992, 329558.016, 249880)
    6 = @(127911, 138648.992, 150315.008)
    7 = @(0, 0, 0)
    8 = @(0, 0, 0)
    9 = @(0, 99358, 56426)
    10 = @(0, 0, 0)
    11 = @(148668.992, 47367, 50073)
    12 = @(56566, 67186, 120034)
    13 = @(0, 0, 0)
    14 = @(0, 0, 0)
    15 = @(0, 0, 0)
    16 = @(547, 561, 50582)
    17 = @(0, 0, 0)
    18 = @(0, 0, 0)
    19 = @(18210, 22453, 23807)
    20 = @(0, 0, 0)
    21 = @(0, 0, 0)
    22 = @(205386, 213106, 215035.008)
    23 = @(12738, 12154, 11005)
    24 = @(1695975.04, 1671147.008, 1652813.952)
    25 = @(0, 0, 0)
    26 = @(2543864.064, 2615521.024, 2532305.92)
    27 = @(595363.968, 588588.992, 544990.0159999999)
    28 = @(46707, 37625, 37581)
    29 = @(142064, 163934, 195474)
    30 = @(45337, 38937, 42642)
    31 = @(37132, 32799, 45536)
    32 = @(0, 0, 0)
    33 = @(68239, 58679, 46752)
    34 = @(255885, 256615.016, 177004.992)
    35 = @(0, 0, 0)
    36 = @(0, 0, 0)
    37 = @(792921.024, 810825.024, 739190.0159999999)
    38 = @(318983.008, 312913.984, 309963.008)
    39 = @(0, 0, 0)
    40 = @(382619.008, 396228, 321576)
    41 = @(72610, 82476, 87171)
    42 = @(0, 0, 0)
    43 = @(18709, 19207, 20480)
    44 = @(0, 0, 0)
    45 = @(0, 0, 0)
    46 = @(673937.024, 696971.008, 699084.032)
    47 = @(481641.9839999999, 519135.9999999999, 549041.92)
    48 = @(40000, 40000, 40000)
    49 = @(33454, 33454, 33454)
    50 = @(0, 0, 0)
    51 = @(407580, 407580, 403599.008)
    52 = @(0, 37494, 71381)
    53 = @(608, 608, 608)
    54 = @(0, 0, 0)
    55 = @(0, 0, 0)
    56 = @(0, 0, 0)
    59 = @(321942.944, 291135.008, 307995.008)
    60 = @(-110100, -115154, -110142)
    61 = @(211843.04, 175980.992, 197852.992)
    62 = @(0, 0, 0)
    63 = @(-38602, -36902, -44202)
    64 = @(0, 0, 0)
    65 = @(18873, 16146, 14632)
    66 = @(-72370, -35850, -38139)
    67 = @(6917, 8693, 7930)
    68 = @(-8610, -21470, -24811)
    69 = @(22138, 13844, 13125)
    70 = @(-30748, -35314, -37936)
    74 = @(118051.008, 106598, 113263)
    75 = @(-31713, -33729, -40721)
    76 = @(-1495, -5622, -3037)
    79 = @(-26676, -29752, -35619)
    80 = @(58167, 37494, 33887)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 38).Value = $vals[0]
    $ws.Cells.Item($row, 39).Value = $vals[1]
    $ws.Cells.Item($row, 40).Value = $vals[2]
}

# --- Section-header / blank rows: keep AL/AM/AN blank like the rest of the row ---
$ws.Range("AK57:AK58").Copy()
$ws.Range("AL57:AN58").PasteSpecial(-4122)
$ws.Range("AK71:AK73").Copy()
$ws.Range("AL71:AN73").PasteSpecial(-4122)
$ws.Range("AK77:AK78").Copy()
$ws.Range("AL77:AN78").PasteSpecial(-4122)
